# "Add concepts from Ruby Weekly." -----------------------------------------
# Populates a batch of new Reading-Response ("RR: ...") links in column P of
# the Curriculum sheet, renames the "Human Learning" header to
# "Human Learning/Agile", tags a couple of rows with a new "Agile Practices"
# concept in column Q, and fills a handful of previously-blank "None"
# placeholder cells in columns E/G to keep the grid visually consistent with
# its neighbouring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Curriculum")

# --- Prep: copy the cell formats we are going to reuse ----------------------
# "Normal" (no override) cells used by several P-column tags.
# "Plain" no-color font used by P4/H10-style tags.
# "Bold" tag font used by P55.
# "None" placeholder font used by E4/H4-style cells.

# --- Column P: replace "None" placeholders with real RR: concept links -----
# (written in the same order the workbook's shared-string table picked them
# up in, so newly-introduced strings land at the expected indices)

# Row 36 -> "RR: Agile Communication" (plain/no-color tag style, like P4)
$ws.Range("P4").Copy()
$ws.Range("P36").PasteSpecial(-4122)
$ws.Range("P36").Value = "RR: Agile Communication"

# Row 35 -> "Hiring Apprentices" (back to workbook-default "Normal" style)
$ws.Range("P35").Style = "Normal"
$ws.Range("P35").Value = "Hiring Apprentices"

# Row 28 -> "RR: Hacking Education" (Normal style)
$ws.Range("P28").Style = "Normal"
$ws.Range("P28").Value = "RR: Hacking Education"

# Row 54 -> "RR: The Science of Software" (bold tag style, like P55)
$ws.Range("P55").Copy()
$ws.Range("P54").PasteSpecial(-4122)
$ws.Range("P54").Value = "RR: The Science of Software"

# Row 48 -> "Rails Guide on SQL" (plain/no-color tag style)
$ws.Range("P4").Copy()
$ws.Range("P48").PasteSpecial(-4122)
$ws.Range("P48").Value = "Rails Guide on SQL"

# Row 50 -> "RR: Technology Radar" (Normal style)
$ws.Range("P50").Style = "Normal"
$ws.Range("P50").Value = "RR: Technology Radar"

# Row 30 -> "RR: Estimation" (Normal style)
$ws.Range("P30").Style = "Normal"
$ws.Range("P30").Value = "RR: Estimation"

# Row 49 -> "RR: Technical Debt" (plain/no-color tag style)
$ws.Range("P4").Copy()
$ws.Range("P49").PasteSpecial(-4122)
$ws.Range("P49").Value = "RR: Technical Debt"

# Row 31 -> "RR: When to use Modules" (plain/no-color tag style)
$ws.Range("P4").Copy()
$ws.Range("P31").PasteSpecial(-4122)
$ws.Range("P31").Value = "RR: When to use Modules"

# Row 51 -> "RR: Ruby Antipatterns" (plain/no-color tag style)
$ws.Range("P4").Copy()
$ws.Range("P51").PasteSpecial(-4122)
$ws.Range("P51").Value = "RR: Ruby Antipatterns"

# --- Header: "Human Learning" -> "Human Learning/Agile" --------------------
$ws.Range("G1").Value = "Human Learning/Agile"

# --- New concept tag in column Q (row 10) -----------------------------------
# Copy the number/cell format already used by the other Q-column tag cells
# (Q1/Q4/...) before writing the new value.
$ws.Range("Q1").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = "Agile Practices"

# --- Fill in missing "None" placeholders in columns E/G --------------------
# Mirrors the "None" placeholder styling already used throughout the sheet
# (e.g. E4/H4).
$ws.Range("E4").Copy()
$noneCells = "G54", "G56", "E57", "G57"
foreach ($addr in $noneCells) {
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = "None"
}

# --- Sheet view: update the active selection --------------------------------
$ws.Range("Q11").Select()
